# Update "Estado de Cuenta" workbook:
#  - Remove the second worker row (EDILBERTO BUITRAGO SUAREZ / 79349507 / 2504)
#  - Update "VALOR MORA" total, and worker/period counters to reflect the
#    remaining single record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the EDILBERTO BUITRAGO SUAREZ data row entirely (row 17). This also
# shifts the signature block (rows 22/23) up to rows 21/22 and keeps their
# merged ranges intact.
$ws.Rows(17).Delete()

# Refresh the summary figures now that only one worker/period remains.
$ws.Range("E11").Value = 6554
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
